$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "68.279.40"
Set-TextCell "E2" "  +1.68%  "

Set-TextCell "D3" "2.508.85"
Set-TextCell "E3" "  +1.62%  "

Set-TextCell "E4" "  +0.04%  "

Set-TextCell "D5" "591.33"
Set-TextCell "E5" "  +1.45%  "

Set-TextCell "D6" "176.86"
Set-TextCell "E6" "  +1.45%  "

Set-TextCell "E8" "  +1.03%  "

Set-TextCell "D9" "2.508.93"
Set-TextCell "E9" "  +1.62%  "

Set-TextCell "D10" "0.142"
Set-TextCell "E10" "  +3.33%  "

Set-TextCell "E11" "  -0.87%  "

Set-TextCell "E12" "  +0.64%  "

Set-TextCell "D13" "0.335"
Set-TextCell "E13" "  +0.95%  "

Set-TextCell "D14" "2.991.38"
Set-TextCell "E14" "  +2.14%  "

Set-TextCell "D15" "25.82"
Set-TextCell "E15" "  +1.67%  "

Set-TextCell "D16" "68.062.92"
Set-TextCell "E16" "  +1.51%  "

Set-TextCell "E17" "  +0.30%  "

Set-TextCell "D18" "2.503.87"
Set-TextCell "E18" "  +1.25%  "

Set-TextCell "D19" "11.00"
Set-TextCell "E19" "  +0.83%  "

Set-TextCell "D20" "7.44"
Set-TextCell "E20" "  -0.62%  "

Set-TextCell "D21" "351.15"
Set-TextCell "E21" "  +0.56%  "

Set-TextCell "E22" "  +3.85%  "

Set-TextCell "D23" "71.32"
Set-TextCell "E23" "  +3.06%  "

Set-TextCell "D24" "1.00"
Set-TextCell "E24" "  -0.01%  "

Set-TextCell "D25" "4.20"
Set-TextCell "E25" "  +0.22%  "

Set-TextCell "D26" "1.72"
Set-TextCell "E26" "  -4.23%  "

Set-TextCell "D27" "9.13"
Set-TextCell "E27" "  +0.14%  "

Set-TextCell "E28" "  +1.76%  "

Set-TextCell "D29" "0.999"
Set-TextCell "E29" "  -0.10%  "

Set-TextCell "D30" "0.0₃0895"
Set-TextCell "E30" "  -0.64%  "

Set-TextCell "D31" "509.80"
Set-TextCell "E31" "  +1.87%  "

Set-TextCell "D32" "7.78"
Set-TextCell "E32" "  +0.62%  "

Set-TextCell "E33" "  +2.08%  "

Set-TextCell "D34" "1.77"
Set-TextCell "E34" "  +1.08%  "

Set-TextCell "E35" "  +0.01%  "

Set-TextCell "E36" "  +0.50%  "

Set-TextCell "D37" "162.00"
Set-TextCell "E37" "  +0.00%  "

Set-TextCell "D38" "18.69"
Set-TextCell "E38" "  +0.05%  "

Set-TextCell "D39" "18.35"
Set-TextCell "E39" "  +1.25%  "

Set-TextCell "D40" "1.32"
Set-TextCell "E40" "  -0.23%  "

Set-TextCell "E41" "  +0.02%  "

Set-TextCell "D42" "1.74"
Set-TextCell "E42" "  +3.24%  "

Set-TextCell "D43" "0.328"
Set-TextCell "E43" "  +0.21%  "

Set-TextCell "D44" "4.81"
Set-TextCell "E44" "  -0.15%  "

Set-TextCell "D45" "2.42"
Set-TextCell "E45" "  +1.93%  "

Set-TextCell "D46" "151.26"
Set-TextCell "E46" "  +6.44%  "

Set-TextCell "D47" "3.56"
Set-TextCell "E47" "  +2.46%  "

Set-TextCell "B48" "BabyDogeCoin"
Set-TextCell "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D48" "0.0₆0259"
Set-TextCell "E48" "  +2.40%  "

Set-TextCell "B49" "ARBITRUM"
Set-TextCell "C49" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D49" "0.518"
Set-TextCell "E49" "  +1.41%  "

Set-TextCell "B50" "Optimism"
Set-TextCell "C50" "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
Set-TextCell "D50" "1.60"
Set-TextCell "E50" "  +1.74%  "

Set-TextCell "B51" "Cronos"
Set-TextCell "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D51" "0.0740"
Set-TextCell "E51" "  +0.27%  "
